{"js": "// Update the title/date paragraph at the top of the document.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2022-11-30 Wednesday\", \"Replace\");\n\n// Update every math-problem cell in the single 20-row x 5-column table, in\n// row-major order, preserving each cell's existing run formatting (the\n// value setter only rewrites the cell Range's text).\nconst newValues = [\n  \"84-15=\", \"94-37=\", \"58+16=\", \"48+29=\", \"7-0=\",\n  \"30+58=\", \"51-33=\", \"84-68=\", \"85-46=\", \"71+11=\",\n  \"52+24=\", \"86-5=\", \"25+24=\", \"30+39=\", \"2+74=\",\n  \"0+62=\", \"10+43=\", \"23+25=\", \"4+92=\", \"9+9=\",\n  \"62-35=\", \"58-38=\", \"60+37=\", \"64+0=\", \"19+43=\",\n  \"21+11=\", \"26+56=\", \"90-17=\", \"17+76=\", \"3+61=\",\n  \"74-11=\", \"74-63=\", \"97-55=\", \"69-33=\", \"67-26=\",\n  \"93-88=\", \"90-57=\", \"10+5=\", \"77-42=\", \"58+30=\",\n  \"1+0=\", \"20+75=\", \"70-57=\", \"13+82=\", \"17-15=\",\n  \"99-99=\", \"8+17=\", \"40+13=\", \"27-27=\", \"53-38=\",\n  \"57-23=\", \"87-34=\", \"89-27=\", \"97-80=\", \"12+49=\",\n  \"24+33=\", \"83-82=\", \"58+37=\", \"31-0=\", \"26+28=\",\n  \"8+38=\", \"38-7=\", \"45+34=\", \"45-9=\", \"68-42=\",\n  \"97-64=\", \"48-21=\", \"29+42=\", \"66-28=\", \"35+9=\",\n  \"26+35=\", \"18+46=\", \"91-56=\", \"61+36=\", \"64-17=\",\n  \"81-38=\", \"32+5=\", \"38+3=\", \"62+2=\", \"25+5=\",\n  \"41+29=\", \"98-94=\", \"59+10=\", \"48+26=\", \"59-57=\",\n  \"49+28=\", \"74-2=\", \"3+1=\", \"18+29=\", \"49+26=\",\n  \"1+53=\", \"9+40=\", \"37+1=\", \"82-80=\", \"29-1=\",\n  \"63-51=\", \"81+3=\", \"77-21=\", \"69-3=\", \"57+39=\"\n];\n\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rowCount = table.rows.items.length;\nconst colCount = 5;\n\nlet i = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    table.getCell(r, c).value = newValues[i];\n    i++;\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the title/date paragraph at the top of the document.\n$titleRange = $d.Paragraphs.Item(1).Range\n$titleRange.Find.Execute(\"2022-11-29 Tuesday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2022-11-30 Wednesday\", 2) | Out-Null\n\n# Auto-generated replacement values for the 20x5 table (row-major order)\n$newValues = @(\n    \"84-15=\",\n    \"94-37=\",\n    \"58+16=\",\n    \"48+29=\",\n    \"7-0=\",\n    \"30+58=\",\n    \"51-33=\",\n    \"84-68=\",\n    \"85-46=\",\n    \"71+11=\",\n    \"52+24=\",\n    \"86-5=\",\n    \"25+24=\",\n    \"30+39=\",\n    \"2+74=\",\n    \"0+62=\",\n    \"10+43=\",\n    \"23+25=\",\n    \"4+92=\",\n    \"9+9=\",\n    \"62-35=\",\n    \"58-38=\",\n    \"60+37=\",\n    \"64+0=\",\n    \"19+43=\",\n    \"21+11=\",\n    \"26+56=\",\n    \"90-17=\",\n    \"17+76=\",\n    \"3+61=\",\n    \"74-11=\",\n    \"74-63=\",\n    \"97-55=\",\n    \"69-33=\",\n    \"67-26=\",\n    \"93-88=\",\n    \"90-57=\",\n    \"10+5=\",\n    \"77-42=\",\n    \"58+30=\",\n    \"1+0=\",\n    \"20+75=\",\n    \"70-57=\",\n    \"13+82=\",\n    \"17-15=\",\n    \"99-99=\",\n    \"8+17=\",\n    \"40+13=\",\n    \"27-27=\",\n    \"53-38=\",\n    \"57-23=\",\n    \"87-34=\",\n    \"89-27=\",\n    \"97-80=\",\n    \"12+49=\",\n    \"24+33=\",\n    \"83-82=\",\n    \"58+37=\",\n    \"31-0=\",\n    \"26+28=\",\n    \"8+38=\",\n    \"38-7=\",\n    \"45+34=\",\n    \"45-9=\",\n    \"68-42=\",\n    \"97-64=\",\n    \"48-21=\",\n    \"29+42=\",\n    \"66-28=\",\n    \"35+9=\",\n    \"26+35=\",\n    \"18+46=\",\n    \"91-56=\",\n    \"61+36=\",\n    \"64-17=\",\n    \"81-38=\",\n    \"32+5=\",\n    \"38+3=\",\n    \"62+2=\",\n    \"25+5=\",\n    \"41+29=\",\n    \"98-94=\",\n    \"59+10=\",\n    \"48+26=\",\n    \"59-57=\",\n    \"49+28=\",\n    \"74-2=\",\n    \"3+1=\",\n    \"18+29=\",\n    \"49+26=\",\n    \"1+53=\",\n    \"9+40=\",\n    \"37+1=\",\n    \"82-80=\",\n    \"29-1=\",\n    \"63-51=\",\n    \"81+3=\",\n    \"77-21=\",\n    \"69-3=\",\n    \"57+39=\"\n)\n\n# Walk the single 20-row x 5-column table in row-major order (matches $newValues)\n# and overwrite each cell's math-problem text, preserving the cell's existing\n# run formatting (font/size) since we only touch the Range.Text of the run.\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $table.Cell($r, $c)\n        $cell.Range.Text = $newValues[$i]\n        $i++\n    }\n}\n"}
